$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer date field: 10/20/2024 -> 10/21/2024 on the slide master and on
#    every slide layout's "Date Placeholder" shape.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster

for ($mi = 1; $mi -le $master.Shapes.Count; $mi++) {
    $msh = $master.Shapes.Item($mi)
    if ($msh.Name -like "Date Placeholder*") {
        $msh.TextFrame.TextRange.Text = "10/21/2024"
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($lj = 1; $lj -le $layout.Shapes.Count; $lj++) {
        $lsh = $layout.Shapes.Item($lj)
        if ($lsh.Name -like "Date Placeholder*") {
            $lsh.TextFrame.TextRange.Text = "10/21/2024"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 2 ("Pipeline" diagram slide): the whole diagram content shifted
#    up slightly (uniform vertical shift). Update each shape's Top.
#    Literal point values below are chosen so that, after round-tripping
#    through the host's single-precision Top/Left storage, they land on
#    the exact target EMU offset recorded in the authoritative XML.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(2)

$s.Shapes.Item(1).Top  = [double]36.55204772949219    # Rectangle: Rounded Corners 37
$s.Shapes.Item(2).Top  = [double]316.9794616699219    # TextBox 38 ("Prompt Instructions")
$s.Shapes.Item(3).Top  = [double]338.2074890136719    # Rectangle: Rounded Corners 39
$s.Shapes.Item(10).Top = [double]5.630393981933594    # TextBox 67
$s.Shapes.Item(11).Top = [double]38.719844818115234   # Rectangle: Rounded Corners 69
$s.Shapes.Item(14).Top = [double]5.630393981933594    # TextBox 85
$s.Shapes.Item(15).Top = [double]444.6144104003906    # Rectangle: Rounded Corners 100
$s.Shapes.Item(16).Top = [double]448.5816650390625    # Table 103
$s.Shapes.Item(18).Top = [double]421.7299499511719    # TextBox 106
$s.Shapes.Item(19).Top = [double]420.1004943847656    # Picture 107

# Connector: Curved 77 - reposition and give it an explicit accent1 @ 64%
# line colour (it previously relied on the style's lnRef for colour).
$conn78 = $s.Shapes.Item(12)
$conn78.Top = [double]202.34205627441406
$conn78.Line.Transparency = 0.36
$conn78.Line.ForeColor.ObjectThemeColor = 5

# Connector: Curved 80 - same treatment.
$conn81 = $s.Shapes.Item(13)
$conn81.Top = [double]202.34197998046875
$conn81.Line.Transparency = 0.36
$conn81.Line.ForeColor.ObjectThemeColor = 5

# Straight Arrow Connector 105 - same treatment.
$conn106 = $s.Shapes.Item(17)
$conn106.Top = [double]410.2372741699219
$conn106.Line.Transparency = 0.36
$conn106.Line.ForeColor.ObjectThemeColor = 5

# ---------------------------------------------------------------------------
# 3) New "Information Context" label textbox, mirroring the existing
#    "Prompt Instructions" label but over the right-hand column.
# ---------------------------------------------------------------------------
$newTb = $s.Shapes.AddTextbox(1, [double]398.2984313964844, [double]316.9794616699219, [double]159.54969787597656, [double]19.387481689453125)
$newTb.Fill.Visible = 0
$newTf = $newTb.TextFrame
$newTf.WordWrap = -1
$newTf.AutoSize = 1
$newTr = $newTf.TextRange
$newTr.Text = "Information Context"
$newTr.ParagraphFormat.Alignment = 2
$newTr.Font.Size = 10
$newTr.Font.Bold = -1
